$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 161
$ws.Range("F5").Value = 450
$ws.Range("F6").Value = 821
$ws.Range("F8").Value = 1185
$ws.Range("F9").Value = 338
$ws.Range("F12").Value = 681
$ws.Range("F14").Value = 504
$ws.Range("F15").Value = 141
$ws.Range("F18").Value = 2914
$ws.Range("F19").Value = 2619
$ws.Range("F24").Value = 225
$ws.Range("F26").Value = 5262
$ws.Range("F28").Value = 980
$ws.Range("F30").Value = 57
$ws.Range("F31").Value = 304
$ws.Range("F32").Value = 1088
$ws.Range("F35").Value = 287

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 1117
$ws.Range("F5").Value = 5
$ws.Range("F24").Value = 314
$ws.Range("F26").Value = 3920
$ws.Range("F28").Value = 6

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 1780
$ws.Range("F5").Value = 2451
$ws.Range("F6").Value = 1035
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = 1319
$ws.Range("F10").Value = 358
$ws.Range("F11").Value = 100

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1780
$ws.Range("F4").Value = 2451
$ws.Range("F6").Value = 1035
$ws.Range("F7").Value = 1319
$ws.Range("F8").Value = 358
$ws.Range("F9").Value = 100
$ws.Range("F10").Value = 161
$ws.Range("F11").Value = 451
$ws.Range("F12").Value = 821
$ws.Range("F14").Value = 1185
$ws.Range("F15").Value = 338
$ws.Range("F17").Value = 681
$ws.Range("F18").Value = 1117
$ws.Range("F19").Value = 1117
$ws.Range("F20").Value = 504
$ws.Range("F23").Value = 2914
$ws.Range("F24").Value = 2619
$ws.Range("F29").Value = 225
$ws.Range("F30").Value = 5262
$ws.Range("F32").Value = 980
$ws.Range("F37").Value = 57
$ws.Range("F38").Value = 304
$ws.Range("F44").Value = 314
$ws.Range("F45").Value = 314
$ws.Range("F47").Value = 1088
$ws.Range("F51").Value = 287

